$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so that
# numeric-looking strings (e.g. "583.73") are NOT auto-converted
# to floating point numbers by Excel, matching the original
# inlineStr text cells from the workbook.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 2 4 "60.321.88"
Set-TextValue 2 5 "  -0.59%  "

Set-TextValue 3 4 "2.604.88"
Set-TextValue 3 5 "  +0.04%  "

Set-TextValue 4 5 "  -0.04%  "

Set-TextValue 5 4 "583.73"
Set-TextValue 5 5 "  +2.36%  "

Set-TextValue 6 4 "142.91"
Set-TextValue 6 5 "  +0.08%  "

Set-TextValue 7 5 "  +0.19%  "

Set-TextValue 8 5 "  -0.89%  "

Set-TextValue 9 4 "6.52"
Set-TextValue 9 5 "  +0.23%  "

Set-TextValue 10 5 "  -1.66%  "

Set-TextValue 11 4 "0.156"
Set-TextValue 11 5 "  +0.21%  "

Set-TextValue 12 4 "0.373"
Set-TextValue 12 5 "  +1.63%  "

Set-TextValue 13 4 "3.064.53"
Set-TextValue 13 5 "  -0.22%  "

Set-TextValue 14 4 "24.65"
Set-TextValue 14 5 "  +4.64%  "

Set-TextValue 15 4 "60.298.17"
Set-TextValue 15 5 "  -0.64%  "

Set-TextValue 16 5 "  +0.15%  "

Set-TextValue 17 4 "2.608.31"
Set-TextValue 17 5 "  -0.43%  "

Set-TextValue 18 4 "11.38"
Set-TextValue 18 5 "  +1.64%  "

Set-TextValue 19 5 "  -1.21%  "

Set-TextValue 20 4 "346.83"
Set-TextValue 20 5 "  -0.16%  "

Set-TextValue 21 5 "  -2.57%  "

Set-TextValue 22 5 "  -0.18%  "

Set-TextValue 23 5 "  +3.46%  "

Set-TextValue 24 4 "63.94"
Set-TextValue 24 5 "  +0.29%  "

Set-TextValue 25 4 "0.996"
Set-TextValue 25 5 "  +0.06%  "

Set-TextValue 26 5 "  -0.23%  "

Set-TextValue 27 4 "7.97"
Set-TextValue 27 5 "  +3.36%  "

Set-TextValue 28 5 "  +2.59%  "

Set-TextValue 29 4 "0.0₃0797"
Set-TextValue 29 5 "  +0.21%  "

Set-TextValue 30 4 "169.35"
Set-TextValue 30 5 "  +4.78%  "

Set-TextValue 31 4 "6.42"
Set-TextValue 31 5 "  +1.98%  "

Set-TextValue 32 5 "  +0.14%  "

Set-TextValue 33 4 "19.43"
Set-TextValue 33 5 "  -0.43%  "

Set-TextValue 34 5 "  +9.41%  "

Set-TextValue 35 4 "4.27"
Set-TextValue 35 5 "  +0.10%  "

Set-TextValue 36 4 "0.994"
Set-TextValue 36 5 "  +3.76%  "

Set-TextValue 37 5 "  +3.22%  "

Set-TextValue 38 4 "316.36"
Set-TextValue 38 5 "  +6.52%  "

Set-TextValue 40 4 "3.89"
Set-TextValue 40 5 "  +2.13%  "

Set-TextValue 41 4 "0.849"

Set-TextValue 42 4 "135.57"
Set-TextValue 42 5 "  -2.50%  "

Set-TextValue 43 4 "0.0995"
Set-TextValue 43 5 "  +1.12%  "

Set-TextValue 44 5 "  +0.34%  "

Set-TextValue 45 5 "  +1.57%  "

Set-TextValue 46 5 "  +0.48%  "

Set-TextValue 47 2 "RenderToken"
Set-TextValue 47 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 47 4 "5.01"
Set-TextValue 47 5 "  +4.27%  "

Set-TextValue 48 2 "Hedera"
Set-TextValue 48 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 48 4 "0.0548"
Set-TextValue 48 5 "  -0.57%  "

Set-TextValue 49 2 "InjectiveProtocol"
Set-TextValue 49 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 49 4 "19.99"
Set-TextValue 49 5 "  +1.70%  "

Set-TextValue 50 2 "VeChain"
Set-TextValue 50 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 50 4 "0.0241"
Set-TextValue 50 5 "  -0.16%  "

Set-TextValue 51 4 "10.74"
Set-TextValue 51 5 "  +0.40%  "
